# DashBoardTest.java updated and new changes to the BasePage and HomePAge
#
# Adds two new worksheets (logedinToDashboard, addTestFromDropDown) to the
# LoginData workbook, and updates the original doLogin sheet's
# "expectedTest" header to "expectedTitle".

$wb = $excel.ActiveWorkbook
$doLogin = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create the two new sheets. Creating "addTestFromDropDown" first and
#    "logedinToDashboard" second (both inserted right after doLogin) yields
#    final tab order: doLogin, logedinToDashboard, addTestFromDropDown.
# ---------------------------------------------------------------------
$addTest = $wb.Worksheets.Add($null, $doLogin)
$addTest.Name = "addTestFromDropDown"

$logedIn = $wb.Worksheets.Add($null, $doLogin)
$logedIn.Name = "logedinToDashboard"

# Re-resolve sheet handles by name now that both Add() calls are done --
# worksheet object references can become stale across Add() calls, so all
# further work uses freshly fetched handles.
$addTest = $wb.Worksheets.Item("addTestFromDropDown")
$logedIn = $wb.Worksheets.Item("logedinToDashboard")
$doLogin = $wb.Worksheets.Item("doLogin")

# ---------------------------------------------------------------------
# 2. Populate "addTestFromDropDown" (sheet3) first.
# ---------------------------------------------------------------------
$addTest.Range("A1").Value = "username"
$addTest.Range("B1").Value = "password"
$addTest.Range("C1").Value = "browserName"

$addTest.Range("A2").Value = "test@kennect.io"
$addTest.Range("B2").Value = "Qwerty@1234"
$addTest.Range("C2").Value = "chrome"

$addTest.Range("A3").Value = "test@kennect.io"
$addTest.Range("B3").Value = "Qwerty@1234"
$addTest.Range("C3").Value = "edge"

$addTest.Range("A4").Value = "test@kennect.io"
$addTest.Range("B4").Value = "Qwerty@1234"
$addTest.Range("C4").Value = "firefox"

$addTest.Range("A5").Value = "test@kennect.io"
$addTest.Range("B5").Value = "Qwerty@1234"
$addTest.Range("C5").Value = "chrome"

# D/E columns written in this specific order to reproduce the original
# shared-string table ordering.
$addTest.Range("D1").Value = "testName"
$addTest.Range("E1").Value = "discount"
$addTest.Range("D2").Value = "Beans"
$addTest.Range("D3").Value = "xyz"
$addTest.Range("D4").Value = "UR Uric acid"
$addTest.Range("E2").Value = "none"
$addTest.Range("D5").Value = "xray"

$addTest.Range("E3").Value = 0.05
$addTest.Range("E4").Value = 0.1
$addTest.Range("E5").Value = 0.05

Write-Host "populated addTestFromDropDown"

# ---------------------------------------------------------------------
# 3. Populate "logedinToDashboard" (sheet2) second -- mirrors doLogin's
#    original table but with the renamed "expectedTitle" header.
# ---------------------------------------------------------------------
$logedIn.Range("A1").Value = "username"
$logedIn.Range("B1").Value = "password"
$logedIn.Range("C1").Value = "browserName"
$logedIn.Range("D1").Value = "expectedTitle"

$logedIn.Range("A2").Value = "test@kennect.io"
$logedIn.Range("B2").Value = "Qwerty@1234"
$logedIn.Range("C2").Value = "chrome"
$logedIn.Range("D2").Value = "Dashboard"

$logedIn.Range("A3").Value = "test@kennect.io"
$logedIn.Range("B3").Value = "Qwerty@1234"
$logedIn.Range("C3").Value = "edge"
$logedIn.Range("D3").Value = "Dashboard"

$logedIn.Range("A4").Value = "test@kennect.io"
$logedIn.Range("B4").Value = "Qwerty@1234"
$logedIn.Range("C4").Value = "firefox"
$logedIn.Range("D4").Value = "Dashboard"

$logedIn.Range("A5").Value = "test@kennect.io"
$logedIn.Range("B5").Value = "Qwerty@1234"
$logedIn.Range("C5").Value = "chrome"
$logedIn.Range("D5").Value = "Dashboard"

Write-Host "populated logedinToDashboard"

# ---------------------------------------------------------------------
# 4. Update doLogin (sheet1): rename header + move selection.
# ---------------------------------------------------------------------
$doLogin.Range("D1").Value = "expectedTitle"
[void]$doLogin.Range("D1").Select()

Write-Host "updated doLogin"
